$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new date value in column A, row 6 (next empty row after existing data)
# Use a leading apostrophe so Excel stores the value as plain text (not an
# auto-converted date serial number), consistent with the existing date
# cells in column A, and without altering the cell's style/format.
$ws.Range("A6").Value = "'2022-08-17"
$ws.Range("A6").Style = "Normal"
